# Fix Training Data Issue (#48)
#
# The "Date" column on the sheet holds a literal (text) date stamp that was
# derived incorrectly from the stats page ("5-17-2013-14" - a mangled mix of
# the game date and the season). It needs to be corrected to the real,
# ISO-style game date "2014-05-17" for every data row, while remaining a
# plain text value (not a real Excel date/number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "5-17-2013-14"
$newValue = "2014-05-17"

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1
$firstCol = $usedRange.Column
$lastCol = $firstCol + $usedRange.Columns.Count - 1

# Locate the "Date" header cell on the header row so the fix keeps working
# even if columns shift around; fall back to column BF (58) if not found.
$headerRow = $ws.Range($ws.Cells.Item($firstRow, $firstCol), $ws.Cells.Item($firstRow, $lastCol))
$dateHeader = $headerRow.Find("Date")
if ($dateHeader -ne $null) {
    $dateCol = $dateHeader.Column
} else {
    $dateCol = 58
}

for ($row = $firstRow + 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $dateCol)
    if ($cell.Value2 -eq $oldValue) {
        # Force text ("@") formatting before writing the replacement so Excel
        # doesn't reinterpret the ISO-looking string as a real date serial,
        # then drop straight back to the default "Normal" style so the cell
        # keeps its original (unstyled) appearance.
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
        $cell.Style = "Normal"
    }
}
